# Refresh the cryptocurrency price/volume snapshot on Sheet1.
#
# Column D ("Price") and column E ("Volume(1h)") hold plain text in the
# source workbook (prices use "." as a thousands separator in several
# rows, e.g. "28.440.00", and volumes are padded percent strings like
# "  +0.52%  "). We must keep every updated cell as TEXT, exactly like the
# upstream data-refresh bot does, rather than let Excel reinterpret a
# value such as "317.48" as a number.
#
# Strategy per cell:
#   - If the new text cannot be mistaken for a plain number (it already
#     contains more than one "." or a "%"/space), a normal `.Value =`
#     assignment is safe -- Excel leaves it as text.
#   - If the new text DOES look like a plain number (single "." or none),
#     assigning it directly would silently convert the cell to a numeric
#     value. To avoid that we first write it as a literal-string formula
#     (`="317.48"`), then Copy / Paste-Special-Values it onto itself. That
#     collapses the formula down to a plain text cell without Excel ever
#     getting a chance to "type" the text as a number, and without
#     touching the cell's number format/style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "28.469.83" },
    @{ Cell = "D3"; Value = "1.828.16" },
    @{ Cell = "E3"; Value = "  +2.02%  " },
    @{ Cell = "D4"; Value = "1.001" },
    @{ Cell = "E4"; Value = "  -0.06%  " },
    @{ Cell = "D5"; Value = "317.48" },
    @{ Cell = "E5"; Value = "  +0.62%  " },
    @{ Cell = "E6"; Value = "  +0.00%  " },
    @{ Cell = "D7"; Value = "0.5338" },
    @{ Cell = "E7"; Value = "  -0.54%  " },
    @{ Cell = "D8"; Value = "0.4054" },
    @{ Cell = "E8"; Value = "  +7.70%  " },
    @{ Cell = "D9"; Value = "0.07633" },
    @{ Cell = "E9"; Value = "  +1.99%  " },
    @{ Cell = "D10"; Value = "41.85" },
    @{ Cell = "E10"; Value = "  +1.09%  " },
    @{ Cell = "E11"; Value = "  +1.23%  " },
    @{ Cell = "D12"; Value = "6.349" },
    @{ Cell = "E12"; Value = "  +4.14%  " },
    @{ Cell = "D13"; Value = "1.002" },
    @{ Cell = "E13"; Value = "  -0.03%  " },
    @{ Cell = "D14"; Value = "20.88" },
    @{ Cell = "E14"; Value = "  +2.09%  " },
    @{ Cell = "D15"; Value = "7.545" },
    @{ Cell = "E15"; Value = "  +3.95%  " },
    @{ Cell = "D16"; Value = "1.819.27" },
    @{ Cell = "E16"; Value = "  +1.88%  " },
    @{ Cell = "D17"; Value = "89.32" },
    @{ Cell = "E17"; Value = "  +0.14%  " },
    @{ Cell = "D18"; Value = "0.00001073" },
    @{ Cell = "E18"; Value = "  +1.56%  " },
    @{ Cell = "D19"; Value = "0.06617" },
    @{ Cell = "E19"; Value = "  +1.85%  " },
    @{ Cell = "E20"; Value = "  +1.23%  " },
    @{ Cell = "E21"; Value = "  -0.08%  " },
    @{ Cell = "D22"; Value = "6.063" },
    @{ Cell = "E22"; Value = "  +2.12%  " },
    @{ Cell = "D23"; Value = "28.476.71" },
    @{ Cell = "E23"; Value = "  +0.48%  " },
    @{ Cell = "D24"; Value = "11.29" },
    @{ Cell = "E24"; Value = "  +1.70%  " },
    @{ Cell = "D25"; Value = "2.150" },
    @{ Cell = "E25"; Value = "  +2.90%  " },
    @{ Cell = "D26"; Value = "2.480" },
    @{ Cell = "E26"; Value = "  +7.89%  " },
    @{ Cell = "D27"; Value = "156.61" },
    @{ Cell = "E27"; Value = "  -1.17%  " },
    @{ Cell = "D29"; Value = "2.028.97" },
    @{ Cell = "E29"; Value = "  +1.91%  " },
    @{ Cell = "D30"; Value = "123.58" },
    @{ Cell = "E30"; Value = "  +1.40%  " },
    @{ Cell = "D31"; Value = "1.122" },
    @{ Cell = "D32"; Value = "0.1095" },
    @{ Cell = "E32"; Value = "  +4.53%  " },
    @{ Cell = "D33"; Value = "5.690" },
    @{ Cell = "E33"; Value = "  +2.81%  " },
    @{ Cell = "D34"; Value = "3.663" },
    @{ Cell = "E34"; Value = "  -0.02%  " },
    @{ Cell = "D35"; Value = "0.07165" },
    @{ Cell = "E35"; Value = "  +11.17%  " },
    @{ Cell = "D36"; Value = "0.2266" },
    @{ Cell = "E36"; Value = "  +0.44%  " },
    @{ Cell = "D37"; Value = "0.02346" },
    @{ Cell = "E37"; Value = "  +3.04%  " },
    @{ Cell = "D38"; Value = "5.234" },
    @{ Cell = "E38"; Value = "  +4.65%  " },
    @{ Cell = "D39"; Value = "8.809" },
    @{ Cell = "E39"; Value = "  +3.54%  " },
    @{ Cell = "D40"; Value = "0.6266" },
    @{ Cell = "E40"; Value = "  +1.78%  " },
    @{ Cell = "D41"; Value = "11.32" },
    @{ Cell = "E41"; Value = "  +2.59%  " },
    @{ Cell = "D42"; Value = "1.183" },
    @{ Cell = "E42"; Value = "  -0.13%  " },
    @{ Cell = "D43"; Value = "1.001" },
    @{ Cell = "E43"; Value = "  -0.04%  " },
    @{ Cell = "D44"; Value = "1.398" },
    @{ Cell = "E44"; Value = "  -2.65%  " },
    @{ Cell = "D45"; Value = "13.38" },
    @{ Cell = "E45"; Value = "  +0.59%  " },
    @{ Cell = "D46"; Value = "3.702" },
    @{ Cell = "E46"; Value = "  +1.01%  " },
    @{ Cell = "D47"; Value = "0.5853" },
    @{ Cell = "E47"; Value = "  +1.56%  " },
    @{ Cell = "D48"; Value = "126.08" },
    @{ Cell = "E48"; Value = "  +0.65%  " },
    @{ Cell = "E49"; Value = "  +3.22%  " },
    @{ Cell = "D50"; Value = "1.198" },
    @{ Cell = "E50"; Value = "  -0.01%  " },
    @{ Cell = "D51"; Value = "0.06901" },
    @{ Cell = "E51"; Value = "  +0.77%  " }
)

foreach ($update in $updates) {
    $cell = $ws.Range($update.Cell)
    $value = $update.Value

    if ($value -match '^-?[0-9]+(\.[0-9]+)?$') {
        # Looks like a plain number -- route through a text-literal formula
        # and flatten it to a value so the cell stays text-typed.
        $escaped = $value.Replace('"', '""')
        $cell.Formula = '="' + $escaped + '"'
        $cell.Copy() | Out-Null
        $cell.PasteSpecial(-4163) | Out-Null
    } else {
        $cell.Value = $value
    }
}
